$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.650.79"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("E2").NumberFormat = "General"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.928.67"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E4").NumberFormat = "General"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.38"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +9.46%  "
$ws.Range("E5").NumberFormat = "General"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.92"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E6").NumberFormat = "General"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E7").NumberFormat = "General"
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E8").NumberFormat = "General"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.723"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("E9").NumberFormat = "General"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000333"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("E11").NumberFormat = "General"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.36"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("E12").NumberFormat = "General"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.563.06"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("E13").NumberFormat = "General"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.30"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.98%  "
$ws.Range("E14").NumberFormat = "General"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.931.30"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.96"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("E16").NumberFormat = "General"
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("E17").NumberFormat = "General"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.21"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.60%  "
$ws.Range("E18").NumberFormat = "General"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.73"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("E19").NumberFormat = "General"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.499.58"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("E20").NumberFormat = "General"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.07"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("E21").NumberFormat = "General"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.36"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.15%  "
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("B23").NumberFormat = "General"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C23").NumberFormat = "General"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.40"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.38%  "
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("B24").NumberFormat = "General"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C24").NumberFormat = "General"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.18"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +16.14%  "
$ws.Range("E24").NumberFormat = "General"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.03"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E25").NumberFormat = "General"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.70"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.69"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.77%  "
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.49"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("E28").NumberFormat = "General"
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "697.26"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E29").NumberFormat = "General"
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.19"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("E30").NumberFormat = "General"
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.127"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("E31").NumberFormat = "General"
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Toncoin"
$ws.Range("B32").NumberFormat = "General"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C32").NumberFormat = "General"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.84"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("E32").NumberFormat = "General"
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "OKB"
$ws.Range("B33").NumberFormat = "General"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C33").NumberFormat = "General"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "69.36"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.71%  "
$ws.Range("E33").NumberFormat = "General"
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.453"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +15.80%  "
$ws.Range("E34").NumberFormat = "General"
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.97"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("E35").NumberFormat = "General"
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.98"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.78%  "
$ws.Range("E36").NumberFormat = "General"
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0840"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("E37").NumberFormat = "General"
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.148"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E38").NumberFormat = "General"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E39").NumberFormat = "General"
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E40").NumberFormat = "General"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0479"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E41").NumberFormat = "General"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.10"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("E42").NumberFormat = "General"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.77"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("E43").NumberFormat = "General"
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Stacks"
$ws.Range("B44").NumberFormat = "General"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C44").NumberFormat = "General"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.23"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +14.09%  "
$ws.Range("E44").NumberFormat = "General"
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("B45").NumberFormat = "General"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C45").NumberFormat = "General"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.96"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.64%  "
$ws.Range("E45").NumberFormat = "General"
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("B46").NumberFormat = "General"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C46").NumberFormat = "General"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.38"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("E46").NumberFormat = "General"
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Stellar"
$ws.Range("B47").NumberFormat = "General"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C47").NumberFormat = "General"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.141"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E47").NumberFormat = "General"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0351"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("E48").NumberFormat = "General"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.30"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("E49").NumberFormat = "General"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.89"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E50").NumberFormat = "General"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.04%  "
$ws.Range("E51").NumberFormat = "General"
$ws.Range("E51").Style = "Normal"
